$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Описание работы) text for rows with changed/truncated content
$ws.Range("E3").Value = ' Мы повысили заработную плату !!!!!!!!!! Логистика – одно из главных направлений в компании «Магнит». Уже более четверти века у нас всё в порядке!   Работа с нами — это:    достойная своевременная зарплата, два раза в месяц. платим гарантированный оклад плюс премию: влияй на свой доход! удобный график: 4/2, 5/5, 6/3 , 4/4 курсы безопасного и экономичного вождения за наш счёт. собственные ремзоны на каждом АТП. вводим в курс дела, даём наставника. бесплатные медосмотры. оформление карт тахографа за счёт компании. выплаты за экономию топлива. акция «Приведи друга» - возможность дополнительного заработка. возможность проявить себя: участвуй в конкурсах профессионального мастерства! полное соблюдение требований ТК РФ    Комфортные условия труда:    удобная спецодежда и защитная обувь. Комп'
$ws.Range("E4").Value = ' Мы повысили заработную плату !!!!!!!!!! Логистика – одно из главных направлений в компании «Магнит». Уже более четверти века у нас всё в порядке! Работа с нами — это:  достойная своевременная зарплата, два раза в месяц. платим гарантированный оклад плюс премию: влияй на свой доход! удобный график: 20/10, 25/5, 28/14, 30/15 курсы безопасного и экономичного вождения за наш счёт. собственные ремзоны на каждом АТП. вводим в курс дела, даём наставника. бесплатные медосмотры. оформление карт тахографа за счёт компании. выплаты за экономию топлива. акция «Приведи друга» - возможность дополнительного заработка. возможность проявить себя: участвуй в конкурсах профессионального мастерства! полное соблюдение требований ТК РФ    Комфортные условия труда:    удобная спецодежда и защитная обувь. Комплекты для зимы: два ви'
$ws.Range("E5").Value = '  Сбер — лучший работодатель России. У нас более 60 компаний экосистемы и 14 тысяч подразделений по всей стране. Присоединяйся к команде малого бизнеса Сбера! На позиции «Менеджера по работе с ключевыми клиентами» ты будешь общаться с предпринимателями и руководителями компаний малого бизнеса и предлагать лучшие решения от Сбера. Эта вакансия для энергичных и общительных людей, готовых вникать в тонкости бизнеса наших клиентов. Мы поможем тебе освоить специальность, сформируем для тебя базу клиентов, предоставим возможности карьерного роста. Тебе предстоит:  развивать бизнес клиентов из закрепленной базы, предлагая продукты и услуги Сбера проводить переговоры с клиентами на территории компаний-партнеров привлекать новых клиентов на об'
$ws.Range("E6").Value = 'Сбер — это 14 тысяч подразделений и лучший работодатель России по версии hh.   Старшие менеджеры по обслуживанию в Сбере – это специалисты, которые проводят операции клиентов в кассе: помогают обменивать деньги, принимают и выдают наличные.   Эта работа подойдет тебе если ты:   · ищешь работу на неполный день или дополнительный заработок к основной работе · хочешь построить карьеру в банковской сфере · ищешь работу с возможностью карьерного перемещения по всей стране · специалист со средним профессиональным, неполным высшим или высшим образованием · хочешь работать в надежной компании · имеешь опыт работы в сфере продаж, обслуживания и консультирования клиентов от полугода.   Тебе предстоит:   · работать в отделении офиса Сбера · проводить кассовые операции: прием, выдача, обмен денег · предлагать продукты и сервисы Сбера · консультировать клиентов по возникающим вопросам.   Работа в Сбере – это:   · стабильный оклад и социальная поддержка сотрудников · возможность выбрать удобный график работы · дополнительный доход ежемесячно · возможностью совмещать работу с другой занятостью и увлечениями · официальное оформление · расширенная пенсионная программа · бесплатная подписка СберПрайм+ и скидки на продукты компаний-партнеров · оплачиваемая мобильная связь: 150 минут и 3 ГБ интернета Сбермобайла · ипотека выгоднее на 4%.   Присоединяйся к команде Сбера!    АтмосфераРабота в Сбере – это:Обучение за счет компании: онлайн курсы в Виртуальной школе Сбера, возможность получить новую квалификациюРеферальная программа для сотрудников: пригласи в команду знакомых профессионалов и получи свое вознаграждение до 100 тыс. рублейКурсы для будущих родителей, материальная поддержка молодых мам, детский отдых и подарки за счет компанииКорпоративная пенсионная программаДелай мир лучше со СберомЗадайте вопрос работодателюОн получит его с откликом на вакансиюГде располагается место работы?Какой график работы?Вакансия открыта?Какая оплата труда?Как с вами связаться?Другой вопросАдресПокрово-МарфиноПоказать на большой карте'
$ws.Range("E9").Value = 'Описание работодателя: ТОГУП "Водное и газовое хозяйство" приглашает на работу секретаря Опыт работы приветствуется. График 5/2 с 08:00 до 17:00. Обязанности:  прием входящей корреспонденции, подготовка документов на подпись руководителю, распределение документооборота, контроль за сроками исполнения поручений руководителя. Требования:  базовые знания офисных программ, ответственность, внимательность, грамотность, дисциплинированность  Задайте вопрос работодателюОн получит его с откликом на вакансиюГде располагается место работы?Какой график работы?Вакансия открыта?Какая оплата труда?Как с вами связаться?Другой вопросАдресТамбов, Студенецкая улица, 3Показать на большой карте'
$ws.Range("E14").Value = 'AstraLab - одна из самых инновационных компаний российского рекламного рынка. Мы первые в России и четвертые в мире научились интегрировать рекламу в фото на базе компьютерного зрения. У нас есть собственная технология AstraVision: компьютерное зрение + обработка языка и текстов. Мы также сотрудничаем с крупнейшим селлером рекламы в стране - ИМХО. С 2018 года мы работаем с клиентами по всему миру и сотрудничаем с крупнейшими сетевыми агентствами (GroupM, BBDO, Dentsu Aegis Network, Publicis Groupe) и брендами (Audi, Lexus, Honda, Pepsi, Nestle, Red Bull, Asus, P&G, Schwarzkopf и др.) Компания Astralab ищет амбициозного Frontend разработчика, который хочет принять участие в создании уникальных и инновационных продуктов. Если вы любите работать с современными технологиями и инструментами, и хотите развиваться в компании с динамичным ростом и высокими стандартами качества, то эта вакансия для вас. Обязанности  Разработка веб-приложений; Оптимизация и модернизация существующих решений; Проведение тестирования; Решение технических вопросов по интеграции продуктов.  Ожидания  Знания JavaScript; Опыт работы с vue.js (vuex, vue router); Знание основных принципов ООП и шаблонов проектирования; Знание HTML5, CSS3, семантическая разметка; Умение работы с веб-шрифтами; Опыт работы с Git.  Плюсом будет  Использо'
$ws.Range("E18").Value = 'У нас полный цикл производства сувенирной продукции. Еженедельно отгружаемся на Wildberries и Ozon. Много ручной работы - сборка, сортировка, накрутка, упаковка. Требуется скорость, аккуратность и внимательность. Есть склад - надо понимать, что, где лежит и в каком количестве. Всё фиксируется в простых эксельках - будет плюсом, если вы шарите в табличках и любите порядок и структуру. Научим формироваь поставки на ВБ и Озон. Пару раз в неделю надо будет сделать отгрузку в траспортную компанию или отъехать по рабочим вопросам. Поэтому наличие своей машины обязательно. Бензин оплачивается. При желании можем научить работать на оборудовании и в принципе погрузить целиком в процесс, если будет интересно развиваться и расти в зарплате. У нас небольшой коллектив, комфортное для работы помещение в районе Ашана. Зимой тепло, летом прохладно, чай/кофе и печеньки. Зарплата каждую неделю без задержек. Оплачиваемая стажировка. График с 9ти до 18ти с перерывом на обед. Среда и воскресенье выходные.Задайте вопрос работодателюОн получит его с откликом на вакансиюГде располагается место работы?Какой график работы?Вакансия открыта?Какая оплата труда?Как с вами связаться?Другой вопросАдресТамбов, Советская улица, 189БПоказать на большой карте'
$ws.Range("E19").Value = 'Обязанности:  Ведение кадрового делопроизводства в полном объеме - 50 чел. Ведение воинского учета; Работа с персональными данными; Ведение штатного расписания, составление, изменение; Подготовка справок, копий трудовых книжек и др. документов с места работы по запросу работников, в т.ч. для работников, выходящих на пенсию по возрасту и досрочно; Ведение табеля учета рабочего времени Подготовка и сдача отчетности: в СФР, сведения в службу занятости населения, в Росстат; Консультация руководства и специалистов учреждения по вопросам, касающихся трудового законодательства; Отчётность с использованием программы Контур-Экстерн. Прием документов и личных заявлений на подпись директору, составление писем, запросов, других документов, подготовка ответов авторам писем по поручению директора. Ведение отдельных служебных поручений директора; Доведение до сведения работников организационных, распорядительных и кадровых документов; Прием больничных листов;  Требования:  Желательно высшее образование, дополнительная подготовка по направлению; Опыт работы от 3-х лет по настоящее время в аналогичной должности по всем кадровым процессам; Отличное знание Трудового кодекса РФ, нормативных актов; Знание изменений и требований к отчетности, предъявляемые в настоящее время; Желание работать и развиваться в своей'

# Clear column F (Навыки) cells that previously held the placeholder value "другое"
$ws.Range("F2").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F18").Value = ""
